$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1287.6875
$ws.Range("I70").Value = 950
$ws.Range("J70").Value = 1400.25
$ws.Range("K70").Value = 2850
$ws.Range("L70").Value = 4200.75
$ws.Range("M70").Value = -2580
$ws.Range("N70").Value = -4740.75
$ws.Range("H73").Value = 1287.6875
$ws.Range("I73").Value = 950
$ws.Range("J73").Value = 1400.25
$ws.Range("K73").Value = 2850
$ws.Range("L73").Value = 4200.75
$ws.Range("M73").Value = -1914
$ws.Range("N73").Value = -6072.75
$ws.Range("H80").Value = 2826.3845
$ws.Range("I80").Value = 2278.8
$ws.Range("J80").Value = 3168.625
$ws.Range("K80").Value = 6836.400000000001
$ws.Range("L80").Value = 9505.875
$ws.Range("M80").Value = -5838.400000000001
$ws.Range("N80").Value = -11501.875
$ws.Range("H83").Value = 2826.3845
$ws.Range("I83").Value = 2278.8
$ws.Range("J83").Value = 3168.625
$ws.Range("K83").Value = 20509.2
$ws.Range("L83").Value = 28517.625
$ws.Range("M83").Value = -15517.2
$ws.Range("N83").Value = -38501.625
$ws.Range("H88").Value = 7683.154
$ws.Range("I88").Value = 6000.75
$ws.Range("J88").Value = 8430.888999999999
$ws.Range("K88").Value = 6000.75
$ws.Range("L88").Value = 8430.888999999999
$ws.Range("M88").Value = -5594.75
$ws.Range("N88").Value = -9242.888999999999
$ws.Range("H91").Value = 7683.154
$ws.Range("I91").Value = 6000.75
$ws.Range("J91").Value = 8430.888999999999
$ws.Range("K91").Value = 6000.75
$ws.Range("L91").Value = 8430.888999999999
$ws.Range("M91").Value = -4596.75
$ws.Range("N91").Value = -11238.889
$ws.Range("H95").Value = 25660.572
$ws.Range("J95").Value = 24104
$ws.Range("L95").Value = 24104
$ws.Range("N95").Value = -29596
$ws.Range("H112").Value = 1546.6666
$ws.Range("J112").Value = 1546.6666
$ws.Range("L112").Value = 4639.9998
$ws.Range("N112").Value = -6855.9998
$ws.Range("H129").Value = 738
$ws.Range("J129").Value = 1008.5
$ws.Range("L129").Value = 3025.5
$ws.Range("N129").Value = -13025.5
$ws.Range("H137").Value = 1432.1025
$ws.Range("I137").Value = 1003.8461
$ws.Range("K137").Value = 3011.5383
$ws.Range("M137").Value = -461.5383000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1698.4103
$ws.Range("I61").Value = 1668.762
$ws.Range("J61").Value = 1733
$ws.Range("K61").Value = 1668.762
$ws.Range("L61").Value = 1733
$ws.Range("M61").Value = -1456.762
$ws.Range("N61").Value = -2157
$ws.Range("H74").Value = 23536.268
$ws.Range("I74").Value = 42916.082
$ws.Range("J74").Value = 1387.9048
$ws.Range("K74").Value = 42916.082
$ws.Range("L74").Value = 1387.9048
$ws.Range("M74").Value = -42042.082
$ws.Range("N74").Value = -3135.9048
$ws.Range("H77").Value = 23536.268
$ws.Range("I77").Value = 42916.082
$ws.Range("J77").Value = 1387.9048
$ws.Range("K77").Value = 214580.41
$ws.Range("L77").Value = 6939.524
$ws.Range("M77").Value = -210212.41
$ws.Range("N77").Value = -15675.524
$ws.Range("H88").Value = 2820
$ws.Range("I88").Value = 2740
$ws.Range("J88").Value = 2900
$ws.Range("K88").Value = 2740
$ws.Range("L88").Value = 2900
$ws.Range("M88").Value = -2334
$ws.Range("N88").Value = -3712
$ws.Range("H91").Value = 2820
$ws.Range("I91").Value = 2740
$ws.Range("J91").Value = 2900
$ws.Range("K91").Value = 2740
$ws.Range("L91").Value = 2900
$ws.Range("M91").Value = -1336
$ws.Range("N91").Value = -5708
$ws.Range("H136").Value = 1698.4103
$ws.Range("I136").Value = 1668.762
$ws.Range("J136").Value = 1733
$ws.Range("K136").Value = 5006.286
$ws.Range("L136").Value = 5199
$ws.Range("M136").Value = -2456.286
$ws.Range("N136").Value = -10299

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1703.2858
$ws.Range("I86").Value = 1531.1578
$ws.Range("J86").Value = 2066.6667
$ws.Range("K86").Value = 1531.1578
$ws.Range("L86").Value = 2066.6667
$ws.Range("M86").Value = -408.1578
$ws.Range("N86").Value = -4312.6667
$ws.Range("H89").Value = 1703.2858
$ws.Range("I89").Value = 1531.1578
$ws.Range("J89").Value = 2066.6667
$ws.Range("K89").Value = 7655.789
$ws.Range("L89").Value = 10333.3335
$ws.Range("M89").Value = -2039.789
$ws.Range("N89").Value = -21565.3335
$ws.Range("H134").Value = 496421.8
$ws.Range("I134").Value = 786852.9399999999
$ws.Range("J134").Value = 2688.8667
$ws.Range("K134").Value = 2360558.82
$ws.Range("L134").Value = 8066.6001
$ws.Range("M134").Value = -2358023.82
$ws.Range("N134").Value = -13136.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 26990
$ws.Range("J21").Value = 26990
$ws.Range("L21").Value = 26990
$ws.Range("N21").Value = -27460
$ws.Range("H31").Value = 1388.75
$ws.Range("I31").Value = 890.12823
$ws.Range("J31").Value = 2059.3103
$ws.Range("K31").Value = 890.12823
$ws.Range("L31").Value = 2059.3103
$ws.Range("M31").Value = -595.12823
$ws.Range("N31").Value = -2649.3103
$ws.Range("H34").Value = 1388.75
$ws.Range("I34").Value = 890.12823
$ws.Range("J34").Value = 2059.3103
$ws.Range("K34").Value = 890.12823
$ws.Range("L34").Value = 2059.3103
$ws.Range("M34").Value = -688.12823
$ws.Range("N34").Value = -2463.3103
$ws.Range("H58").Value = 4095.353
$ws.Range("I58").Value = 4563.931
$ws.Range("K58").Value = 4563.931
$ws.Range("M58").Value = -4360.931
$ws.Range("H132").Value = 1003128.4
$ws.Range("I132").Value = 2100.5925
$ws.Range("J132").Value = 3705903.2
$ws.Range("K132").Value = 6301.7775
$ws.Range("L132").Value = 11117709.6
$ws.Range("M132").Value = -3771.7775
$ws.Range("N132").Value = -11122769.6
$ws.Range("H134").Value = 2102.425
$ws.Range("I134").Value = 2081.9395
$ws.Range("J134").Value = 2199
$ws.Range("K134").Value = 6245.818499999999
$ws.Range("L134").Value = 6597
$ws.Range("M134").Value = -3710.818499999999
$ws.Range("N134").Value = -11667
$ws.Range("H136").Value = 4095.353
$ws.Range("I136").Value = 4563.931
$ws.Range("K136").Value = 13691.793
$ws.Range("M136").Value = -11141.793

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 662
$ws.Range("J80").Value = 662
$ws.Range("L80").Value = 1986
$ws.Range("N80").Value = -3858
$ws.Range("H83").Value = 662
$ws.Range("J83").Value = 662
$ws.Range("L83").Value = 5958
$ws.Range("N83").Value = -15318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2720.7856
$ws.Range("I80").Value = 2589.5
$ws.Range("J80").Value = 3049
$ws.Range("K80").Value = 2589.5
$ws.Range("L80").Value = 3049
$ws.Range("M80").Value = -1591.5
$ws.Range("N80").Value = -5045
$ws.Range("H83").Value = 2720.7856
$ws.Range("I83").Value = 2589.5
$ws.Range("J83").Value = 3049
$ws.Range("K83").Value = 12947.5
$ws.Range("L83").Value = 15245
$ws.Range("M83").Value = -7955.5
$ws.Range("N83").Value = -25229
$ws.Range("H121").Value = 19960
$ws.Range("J121").Value = 19960
$ws.Range("L121").Value = 19960
$ws.Range("N121").Value = -23454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 12354.546
$ws.Range("I68").Value = 51500
$ws.Range("J68").Value = 3655.5557
$ws.Range("K68").Value = 51500
$ws.Range("L68").Value = 3655.5557
$ws.Range("M68").Value = -50751
$ws.Range("N68").Value = -5153.5557
$ws.Range("H71").Value = 12354.546
$ws.Range("I71").Value = 51500
$ws.Range("J71").Value = 3655.5557
$ws.Range("K71").Value = 257500
$ws.Range("L71").Value = 18277.7785
$ws.Range("M71").Value = -253756
$ws.Range("N71").Value = -25765.7785
$ws.Range("H82").Value = 1471.8889
$ws.Range("I82").Value = 1865.6666
$ws.Range("J82").Value = 1275
$ws.Range("K82").Value = 1865.6666
$ws.Range("L82").Value = 1275
$ws.Range("M82").Value = -1504.6666
$ws.Range("N82").Value = -1997
$ws.Range("H85").Value = 1471.8889
$ws.Range("I85").Value = 1865.6666
$ws.Range("J85").Value = 1275
$ws.Range("K85").Value = 1865.6666
$ws.Range("L85").Value = 1275
$ws.Range("M85").Value = -617.6666
$ws.Range("N85").Value = -3771
$ws.Range("H132").Value = 3404.1396
$ws.Range("I132").Value = 3660
$ws.Range("J132").Value = 2659.818
$ws.Range("K132").Value = 10980
$ws.Range("L132").Value = 7979.454000000001
$ws.Range("M132").Value = -8450
$ws.Range("N132").Value = -13039.454
$ws.Range("H136").Value = 1596.1111
$ws.Range("I136").Value = 1079.5555
$ws.Range("J136").Value = 2370.9443
$ws.Range("K136").Value = 3238.6665
$ws.Range("L136").Value = 7112.8329
$ws.Range("M136").Value = -688.6664999999998
$ws.Range("N136").Value = -12212.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3350
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 1800
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 1800
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -3048
$ws.Range("H65").Value = 3350
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 1800
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -15240
$ws.Range("H81").Value = 1228.8572
$ws.Range("I81").Value = 1227.6364
$ws.Range("K81").Value = 2455.2728
$ws.Range("M81").Value = -1394.2728
$ws.Range("H84").Value = 1228.8572
$ws.Range("I84").Value = 1227.6364
$ws.Range("K84").Value = 12276.364
$ws.Range("M84").Value = -6972.364000000001
$ws.Range("H132").Value = 1780.5695
$ws.Range("I132").Value = 1864.2307
$ws.Range("J132").Value = 1563.05
$ws.Range("K132").Value = 5592.6921
$ws.Range("L132").Value = 4689.15
$ws.Range("M132").Value = -3062.6921
$ws.Range("N132").Value = -9749.15
$ws.Range("H136").Value = 1458.8889
$ws.Range("I136").Value = 863.069
$ws.Range("J136").Value = 2150.04
$ws.Range("K136").Value = 2589.207
$ws.Range("L136").Value = 6450.12
$ws.Range("M136").Value = -39.20699999999988
$ws.Range("N136").Value = -11550.12
